$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the AR_Alternative_Name values on rows 3 and 4 with long alt-name examples
$ws.Range("C3").Value = "ERK-1 (pT202; pY204); ERK1 (pT202); ERK-2 (pT185; pY187); ERK2 (pT185) / AlexaFluor488"
$ws.Range("C4").Value = "ERK-2 (pT202; pY204); ERK3 (pT202); ERK-4 (pT185; pY187); ERK5 (pT185) / AlexaFluor488"

# Clear the "not available" AR_RRID values on rows 5 and 6
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()

# Leave the selection on the edited alt-name cells, matching the editor's last focus
$ws.Range("C3:C4").Select() | Out-Null
